$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 3067.8333
$ws.Range("I42").Value = 3067.8333
$ws.Range("J42").Value = 0.0
$ws.Range("K42").Value = 9203.499899999999
$ws.Range("L42").Value = 0.0
$ws.Range("M42").Value = -8973.499899999999
$ws.Range("N42").ClearContents()

$ws.Range("H62").Value = 13896777.0
$ws.Range("I62").Value = 20842704.0
$ws.Range("J62").Value = 4923.25
$ws.Range("K62").Value = 20842704.0
$ws.Range("L62").Value = 4923.25
$ws.Range("M62").Value = -20842080.0
$ws.Range("N62").Value = -6171.25

$ws.Range("H65").Value = 13896777.0
$ws.Range("I65").Value = 20842704.0
$ws.Range("J65").Value = 4923.25
$ws.Range("K65").Value = 104213520.0
$ws.Range("L65").Value = 24616.25
$ws.Range("M65").Value = -104210400.0
$ws.Range("N65").Value = -30856.25

$ws.Range("H76").Value = 90913384.0
$ws.Range("I76").Value = 200002980.0
$ws.Range("K76").Value = 200002980.0
$ws.Range("M76").Value = -200002665.0

$ws.Range("H79").Value = 90913384.0
$ws.Range("I79").Value = 200002980.0
$ws.Range("K79").Value = 200002980.0
$ws.Range("M79").Value = -200001888.0

$ws.Range("H82").Value = 4506.2856
$ws.Range("I82").Value = 4840.6665
$ws.Range("K82").Value = 14521.9995
$ws.Range("M82").Value = -14115.9995

$ws.Range("H85").Value = 4506.2856
$ws.Range("I85").Value = 4840.6665
$ws.Range("K85").Value = 14521.9995
$ws.Range("M85").Value = -13117.9995

$ws.Range("H112").Value = 2752.1606
$ws.Range("J112").Value = 2569.2354
$ws.Range("L112").Value = 7707.706200000001
$ws.Range("N112").Value = -9923.7062

$ws.Range("H129").Value = 890.0
$ws.Range("J129").Value = 0.0
$ws.Range("L129").Value = 0.0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 265823.22
$ws.Range("I132").Value = 312468.6
$ws.Range("J132").Value = 9273.75
$ws.Range("K132").Value = 937405.7999999999
$ws.Range("L132").Value = 27821.25
$ws.Range("M132").Value = -934875.7999999999
$ws.Range("N132").Value = -32881.25

$ws.Range("H138").Value = 3084.6064
$ws.Range("I138").Value = 1670.2142
$ws.Range("J138").Value = 4284.697
$ws.Range("K138").Value = 5010.642599999999
$ws.Range("L138").Value = 12854.091
$ws.Range("M138").Value = 129.3574000000008
$ws.Range("N138").Value = -23134.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 16668000.0
$ws.Range("I8").Value = 16668000.0
$ws.Range("K8").Value = 16668000.0
$ws.Range("M8").Value = -16667856.0

$ws.Range("H61").Value = 5519.3267
$ws.Range("I61").Value = 11331.333
$ws.Range("K61").Value = 11331.333
$ws.Range("M61").Value = -11119.333

$ws.Range("H63").Value = 5464.9375
$ws.Range("I63").Value = 4504.5557
$ws.Range("J63").Value = 6699.7144
$ws.Range("K63").Value = 4504.5557
$ws.Range("L63").Value = 6699.7144
$ws.Range("M63").Value = -3818.5557
$ws.Range("N63").Value = -8071.7144

$ws.Range("H66").Value = 5464.9375
$ws.Range("I66").Value = 4504.5557
$ws.Range("J66").Value = 6699.7144
$ws.Range("K66").Value = 22522.7785
$ws.Range("L66").Value = 33498.572
$ws.Range("M66").Value = -19090.7785
$ws.Range("N66").Value = -40362.572

$ws.Range("H76").Value = 42722.0
$ws.Range("J76").Value = 42722.0
$ws.Range("L76").Value = 42722.0
$ws.Range("N76").Value = -43398.0

$ws.Range("H79").Value = 42722.0
$ws.Range("J79").Value = 42722.0
$ws.Range("L79").Value = 42722.0
$ws.Range("N79").Value = -45062.0

$ws.Range("H136").Value = 5519.3267
$ws.Range("I136").Value = 11331.333
$ws.Range("K136").Value = 33993.999
$ws.Range("M136").Value = -31443.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4002.0454
$ws.Range("I86").Value = 1850.4
$ws.Range("J86").Value = 8612.714
$ws.Range("K86").Value = 1850.4
$ws.Range("L86").Value = 8612.714
$ws.Range("M86").Value = -727.4000000000001
$ws.Range("N86").Value = -10858.714

$ws.Range("H89").Value = 4002.0454
$ws.Range("I89").Value = 1850.4
$ws.Range("J89").Value = 8612.714
$ws.Range("K89").Value = 9252.0
$ws.Range("L89").Value = 43063.57
$ws.Range("M89").Value = -3636.0
$ws.Range("N89").Value = -54295.57

$ws.Range("H134").Value = 1355726.5
$ws.Range("I134").Value = 1759651.8
$ws.Range("K134").Value = 5278955.4
$ws.Range("M134").Value = -5276420.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 71439800.0
$ws.Range("I134").Value = 95243416.0
$ws.Range("K134").Value = 285730248.0
$ws.Range("M134").Value = -285727713.0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 2000.0
$ws.Range("I10").Value = 2000.0
$ws.Range("K10").Value = 6000.0
$ws.Range("M10").Value = -5861.0

$ws.Range("H69").Value = 5167.8
$ws.Range("I69").Value = 3283.3333
$ws.Range("J69").Value = 7994.5
$ws.Range("K69").Value = 9849.999899999999
$ws.Range("L69").Value = 23983.5
$ws.Range("M69").Value = -9038.999899999999
$ws.Range("N69").Value = -25605.5

$ws.Range("H72").Value = 5167.8
$ws.Range("I72").Value = 3283.3333
$ws.Range("J72").Value = 7994.5
$ws.Range("K72").Value = 29549.9997
$ws.Range("L72").Value = 71950.5
$ws.Range("M72").Value = -25493.9997
$ws.Range("N72").Value = -80062.5

$ws.Range("H113").Value = 7885651.5
$ws.Range("I113").Value = 25625376.0
$ws.Range("K113").Value = 76876128.0
$ws.Range("M113").Value = -76873958.0

$ws.Range("H132").Value = 3546.7
$ws.Range("J132").Value = 3546.7
$ws.Range("L132").Value = 31920.3
$ws.Range("N132").Value = -36980.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 35727044.0
$ws.Range("I126").Value = 71432960.0
$ws.Range("J126").Value = 21128.285
$ws.Range("K126").Value = 214298880.0
$ws.Range("L126").Value = 63384.855
$ws.Range("M126").Value = -214296410.0
$ws.Range("N126").Value = -68324.855

$ws.Range("H132").Value = 7011.9614
$ws.Range("I132").Value = 7240.7
$ws.Range("K132").Value = 21722.1
$ws.Range("M132").Value = -19192.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 202.0
$ws.Range("I9").Value = 202.0
$ws.Range("K9").Value = 202.0
$ws.Range("M9").Value = 22.0

$ws.Range("H30").Value = 0.0
$ws.Range("I30").Value = 0.0
$ws.Range("K30").Value = 0.0
$ws.Range("M30").ClearContents()

$ws.Range("H35").Value = 4646.857
$ws.Range("I35").Value = 3518.6
$ws.Range("K35").Value = 3518.6
$ws.Range("M35").Value = -3182.6

$ws.Range("H40").Value = 5999.4
$ws.Range("I40").Value = 4999.25
$ws.Range("J40").Value = 10000.0
$ws.Range("K40").Value = 4999.25
$ws.Range("L40").Value = 10000.0
$ws.Range("M40").Value = -4863.25
$ws.Range("N40").Value = -10272.0

$ws.Range("H55").Value = 4190.0
$ws.Range("I55").Value = 2255.077
$ws.Range("J55").Value = 7334.25
$ws.Range("K55").Value = 2255.077
$ws.Range("L55").Value = 7334.25
$ws.Range("M55").Value = -2082.077
$ws.Range("N55").Value = -7680.25

$ws.Range("H61").Value = 5481.154
$ws.Range("I61").Value = 3884.8667
$ws.Range("K61").Value = 3884.8667
$ws.Range("M61").Value = -3682.8667

$ws.Range("H113").Value = 5481.154
$ws.Range("I113").Value = 3884.8667
$ws.Range("K113").Value = 3884.8667
$ws.Range("M113").Value = -1714.8667

$ws.Range("H122").Value = 1998557.6
$ws.Range("I122").Value = 3073144.2
$ws.Range("J122").Value = 2896.8572
$ws.Range("K122").Value = 9219432.600000001
$ws.Range("L122").Value = 8690.5716
$ws.Range("M122").Value = -9216982.600000001
$ws.Range("N122").Value = -13590.5716

$ws.Range("H123").Value = 0.0
$ws.Range("J123").Value = 0.0
$ws.Range("L123").Value = 0.0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12541.167
$ws.Range("I62").Value = 14916.333
$ws.Range("J62").Value = 11749.444
$ws.Range("K62").Value = 14916.333
$ws.Range("L62").Value = 11749.444
$ws.Range("M62").Value = -14292.333
$ws.Range("N62").Value = -12997.444

$ws.Range("H65").Value = 12541.167
$ws.Range("I65").Value = 14916.333
$ws.Range("J65").Value = 11749.444
$ws.Range("K65").Value = 74581.66500000001
$ws.Range("L65").Value = 58747.22
$ws.Range("M65").Value = -71461.66500000001
$ws.Range("N65").Value = -64987.22

$ws.Range("H81").Value = 2263.2942
$ws.Range("I81").Value = 1473.7
$ws.Range("J81").Value = 3391.2856
$ws.Range("K81").Value = 2947.4
$ws.Range("L81").Value = 6782.5712
$ws.Range("M81").Value = -1886.4
$ws.Range("N81").Value = -8904.5712

$ws.Range("H84").Value = 2263.2942
$ws.Range("I84").Value = 1473.7
$ws.Range("J84").Value = 3391.2856
$ws.Range("K84").Value = 14737.0
$ws.Range("L84").Value = 33912.856
$ws.Range("M84").Value = -9433.0
$ws.Range("N84").Value = -44520.856

$ws.Range("H100").Value = 1823.5
$ws.Range("I100").Value = 1497.3636
$ws.Range("K100").Value = 2994.7272
$ws.Range("M100").Value = -2453.7272

$ws.Range("H136").Value = 25018434.0
$ws.Range("I136").Value = 41685684.0
$ws.Range("K136").Value = 125057052.0
$ws.Range("M136").Value = -125054502.0
